$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume/Change %) updates for the
# crypto ranking refresh. Values that look numeric (e.g. "403.91")
# are forced to Text via NumberFormat so Excel keeps them as literal
# strings (matching the original inlineStr cell content) instead of
# silently converting them to floating point numbers; the format is
# cleared again right after so no stray formatting is left behind.

$ws.Range('D2').Value = '61.536.36'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '3.395.94'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  +0.06%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '403.91'
$cell.ClearFormats()
$ws.Range('E5').Value = '  -0.74%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '130.65'
$cell.ClearFormats()
$ws.Range('E6').Value = '  +0.96%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.589'
$cell.ClearFormats()
$ws.Range('E7').Value = '  -2.39%  '
$ws.Range('E8').Value = '  +0.01%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.679'
$cell.ClearFormats()
$ws.Range('E9').Value = '  -0.22%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.130'
$cell.ClearFormats()
$ws.Range('E10').Value = '  +1.83%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '41.56'
$cell.ClearFormats()
$ws.Range('E11').Value = '  -2.21%  '
$ws.Range('E12').Value = '  -0.76%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '19.72'
$cell.ClearFormats()
$ws.Range('E13').Value = '  -0.67%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '8.33'
$cell.ClearFormats()
$ws.Range('E14').Value = '  -3.59%  '
$ws.Range('D15').Value = '3.394.09'
$ws.Range('E15').Value = '  -1.99%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '11.58'
$cell.ClearFormats()
$ws.Range('E16').Value = '  +1.70%  '
$ws.Range('D17').Value = '61.448.99'
$ws.Range('E17').Value = '  -0.79%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '1.01'
$cell.ClearFormats()
$ws.Range('E18').Value = '  -2.22%  '
$ws.Range('E19').Value = '  +4.54%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '3.16'
$cell.ClearFormats()
$ws.Range('E20').Value = '  -3.31%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '83.03'
$cell.ClearFormats()
$ws.Range('E21').Value = '  -0.43%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '310.83'
$cell.ClearFormats()
$ws.Range('E22').Value = '  +0.39%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '12.69'
$cell.ClearFormats()
$ws.Range('E23').Value = '  -2.56%  '
$ws.Range('E24').Value = '  -1.18%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '4.76'
$cell.ClearFormats()
$ws.Range('E25').Value = '  +7.68%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '8.15'
$cell.ClearFormats()
$ws.Range('E26').Value = '  +8.11%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '29.39'
$cell.ClearFormats()
$ws.Range('E27').Value = '  -1.41%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '8.03'
$cell.ClearFormats()
$ws.Range('E28').Value = '  -6.54%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '2.71'
$cell.ClearFormats()
$ws.Range('E29').Value = '  +6.02%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '44.04'
$cell.ClearFormats()
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('E31').Value = '  -1.90%  '
$ws.Range('E32').Value = '  -2.15%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '11.27'
$cell.ClearFormats()
$ws.Range('E33').Value = '  -3.80%  '
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('E37').Value = '  +0.15%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '2.97'
$cell.ClearFormats()
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('E39').Value = '  -3.81%  '
$ws.Range('E40').Value = '  +10.02%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '139.97'
$cell.ClearFormats()
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('E43').Value = '  -0.88%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '3.92'
$cell.ClearFormats()
$ws.Range('E44').Value = '  +0.19%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '16.69'
$cell.ClearFormats()
$ws.Range('E45').Value = '  -2.56%  '
$ws.Range('E46').Value = '  -1.58%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '21.02'
$cell.ClearFormats()
$ws.Range('E47').Value = '  -4.03%  '
$ws.Range('D48').Value = '2.091.26'
$ws.Range('E48').Value = '  -2.84%  '
$ws.Range('E49').Value = '  -2.28%  '
$ws.Range('E50').Value = '  +2.06%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.68'
$cell.ClearFormats()
$ws.Range('E51').Value = '  +12.39%  '
